$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "27.553.04"
$ws.Cells.Item(2, 5).Value = "  -1.24%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.578.27"
$ws.Cells.Item(3, 5).Value = "  -3.32%  "
$ws.Cells.Item(4, 5).Value = "  +0.24%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "205.88"
$ws.Cells.Item(5, 5).Value = "  -2.61%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.502"
$ws.Cells.Item(6, 5).Value = "  -3.25%  "
$ws.Cells.Item(7, 5).Value = "  +0.28%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "22.09"
$ws.Cells.Item(8, 5).Value = "  -5.88%  "
$ws.Cells.Item(9, 5).Value = "  -2.25%  "
$ws.Cells.Item(10, 5).Value = "  -3.78%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0864"
$ws.Cells.Item(11, 5).Value = "  -2.11%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.801.87"
$ws.Cells.Item(12, 5).Value = "  -3.34%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.590.61"
$ws.Cells.Item(13, 5).Value = "  -2.63%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.84"
$ws.Cells.Item(14, 5).Value = "  -4.77%  "
$ws.Cells.Item(15, 5).Value = "  -7.07%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "27.529.95"
$ws.Cells.Item(16, 5).Value = "  -1.36%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "62.67"
$ws.Cells.Item(17, 5).Value = "  -4.23%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "216.97"
$ws.Cells.Item(18, 5).Value = "  -5.25%  "
$ws.Cells.Item(19, 5).Value = "  -4.07%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.30"
$ws.Cells.Item(20, 5).Value = "  -4.96%  "
$ws.Cells.Item(21, 5).Value = "  +0.36%  "
$ws.Cells.Item(22, 5).Value = "  -4.88%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.41"
$ws.Cells.Item(23, 5).Value = "  -6.44%  "
$ws.Cells.Item(24, 5).Value = "  -4.39%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "153.30"
$ws.Cells.Item(25, 5).Value = "  -1.39%  "
$ws.Cells.Item(26, 5).Value = "  +0.26%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "6.67"
$ws.Cells.Item(27, 5).Value = "  -3.08%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "15.00"
$ws.Cells.Item(28, 5).Value = "  -3.42%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.106"
$ws.Cells.Item(29, 5).Value = "  -4.83%  "
$ws.Cells.Item(30, 5).Value = "  -2.42%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.0461"
$ws.Cells.Item(31, 5).Value = "  -4.07%  "
$ws.Cells.Item(32, 5).Value = "  -5.41%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.364.58"
$ws.Cells.Item(33, 5).Value = "  -1.95%  "
$ws.Cells.Item(34, 5).Value = "  -5.74%  "
$ws.Cells.Item(35, 5).Value = "  -5.56%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.962"
$ws.Cells.Item(36, 5).Value = "  -5.34%  "
$ws.Cells.Item(37, 5).Value = "  -1.48%  "
$ws.Cells.Item(38, 5).Value = "  -4.36%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.534"
$ws.Cells.Item(39, 5).Value = "  -4.34%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.812"
$ws.Cells.Item(40, 5).Value = "  -4.41%  "
$ws.Cells.Item(41, 5).Value = "  +0.29%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.974"
$ws.Cells.Item(42, 5).Value = "  -4.02%  "
$ws.Cells.Item(44, 2).Value = "Aave"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "63.20"
$ws.Cells.Item(44, 5).Value = "  -3.92%  "
$ws.Cells.Item(45, 2).Value = "FraxShare"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "5.25"
$ws.Cells.Item(45, 5).Value = "  -3.28%  "
$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.75"
$ws.Cells.Item(46, 5).Value = "  -3.91%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.712.13"
$ws.Cells.Item(47, 5).Value = "  -3.47%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "87.17"
$ws.Cells.Item(48, 5).Value = "  -1.80%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0₆01000"
$ws.Cells.Item(49, 5).Value = "  -3.23%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0966"
$ws.Cells.Item(50, 5).Value = "  -5.16%  "
$ws.Cells.Item(51, 5).Value = "  -1.73%  "
